$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.414509557120016
$ws.Range("C2").Value = 1.320195973509991
$ws.Range("D2").Value = 3.99283691009699
$ws.Range("E2").Value = 1.346585525682575
$ws.Range("F2").Value = 1.414055654904187
$ws.Range("G2").Value = 1.378958854892205
$ws.Range("H2").Value = 1.346982622402224

$ws.Range("B3").Value = 1.425213289967862
$ws.Range("C3").Value = 1.325315907228027
$ws.Range("D3").Value = 3.23283331388647
$ws.Range("E3").Value = 1.347817107110185
$ws.Range("F3").Value = 1.418476688626426
$ws.Range("G3").Value = 1.372112685789044
$ws.Range("H3").Value = 1.347828200706768

$ws.Range("B4").Value = 1.390138974964715
$ws.Range("C4").Value = 1.36634262392983
$ws.Range("D4").Value = 3.743511689471815
$ws.Range("E4").Value = 1.342443321396813
$ws.Range("F4").Value = 1.393836101180278
$ws.Range("G4").Value = 1.383240013895821
$ws.Range("H4").Value = 1.342718102337791

$ws.Range("B5").Value = 1.423685679465119
$ws.Range("C5").Value = 1.409259220543114
$ws.Range("D5").Value = 2.286912720564253
$ws.Range("E5").Value = 1.350032931609149
$ws.Range("F5").Value = 1.414522073045725
$ws.Range("G5").Value = 1.399850874782828
$ws.Range("H5").Value = 1.349992015156981

$ws.Range("B6").Value = 1.428914950493217
$ws.Range("C6").Value = 1.409908246290728
$ws.Range("D6").Value = 0.5354154108338075
$ws.Range("E6").Value = 1.352104465373867
$ws.Range("F6").Value = 1.419532925750195
$ws.Range("G6").Value = 1.384591965316204
$ws.Range("H6").Value = 1.351843737668598

$ws.Range("B7").Value = 1.413520031540223
$ws.Range("C7").Value = 1.414435199910633
$ws.Range("D7").Value = 0.9101152687398403
$ws.Range("E7").Value = 1.3502050541849
$ws.Range("F7").Value = 1.404978498878209
$ws.Range("G7").Value = 1.390972419996553
$ws.Range("H7").Value = 1.349955890124776

$ws.Range("B8").Value = 1.419552146857844
$ws.Range("C8").Value = 1.286428854980599
$ws.Range("D8").Value = 1.149297017754302
$ws.Range("E8").Value = 1.34067164994444
$ws.Range("F8").Value = 1.409312192220709
$ws.Range("G8").Value = 1.307684694103238
$ws.Range("H8").Value = 1.340676964773855

$ws.Range("B9").Value = 1.426263507133462
$ws.Range("C9").Value = 1.366023141857575
$ws.Range("D9").Value = 0.8337038804962449
$ws.Range("E9").Value = 1.344961376309288
$ws.Range("F9").Value = 1.419164020447879
$ws.Range("G9").Value = 1.359680543457173
$ws.Range("H9").Value = 1.345344167322881

$ws.Range("B10").Value = 1.225441711943755
$ws.Range("C10").Value = 1.409810408930862
$ws.Range("D10").Value = 2.973400820959878
$ws.Range("E10").Value = 1.338037341956309
$ws.Range("F10").Value = 1.27224744976302
$ws.Range("G10").Value = 1.403891052508206
$ws.Range("H10").Value = 1.339021721330763

$ws.Range("B11").Value = 1.183188477810487
$ws.Range("C11").Value = 1.414360880002128
$ws.Range("D11").Value = 2.583480844780979
$ws.Range("E11").Value = 1.339872308311112
$ws.Range("F11").Value = 1.232365762670961
$ws.Range("G11").Value = 1.400737064552017
$ws.Range("H11").Value = 1.340293742975246

$ws.Range("B12").Value = 0.9765449706184237
$ws.Range("C12").Value = 1.405296775565357
$ws.Range("D12").Value = 1.194841888868514
$ws.Range("E12").Value = 1.324747553747766
$ws.Range("F12").Value = 0.9993215601058449
$ws.Range("G12").Value = 1.373175483619289
$ws.Range("H12").Value = 1.324371252076056

$ws.Range("B13").Value = 1.21168381432692
$ws.Range("C13").Value = 1.408286739559427
$ws.Range("D13").Value = 2.47231477791947
$ws.Range("E13").Value = 1.335620606451476
$ws.Range("F13").Value = 1.24896039519903
$ws.Range("G13").Value = 1.392783088423188
$ws.Range("H13").Value = 1.336020726629191

